$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 44018
$ws.Cells.Item(2, 3).Value = 23242
$ws.Cells.Item(2, 5).Value = 2721

# Row 3
$ws.Cells.Item(3, 3).Value = 15102
$ws.Cells.Item(3, 4).Value = 132

# Row 4
$ws.Cells.Item(4, 2).Value = 44018
$ws.Cells.Item(4, 3).Value = 214061
$ws.Cells.Item(4, 4).Value = 18596
$ws.Cells.Item(4, 5).Value = 33265
$ws.Cells.Item(4, 6).Value = 5199
$ws.Cells.Item(4, 7).Value = 30.16
$ws.Cells.Item(4, 11).Value = 110283
$ws.Cells.Item(4, 12).Value = 17048

# Row 6
$ws.Cells.Item(6, 2).Value = 44018
$ws.Cells.Item(6, 3).Value = 52155
$ws.Cells.Item(6, 4).Value = 653
$ws.Cells.Item(6, 5).Value = 10640
$ws.Cells.Item(6, 6).Value = 231
$ws.Cells.Item(6, 7).Value = 20.4
$ws.Cells.Item(6, 8).Value = 35.38

# Row 7
$ws.Cells.Item(7, 2).Value = 44018
$ws.Cells.Item(7, 3).Value = 25469
$ws.Cells.Item(7, 4).Value = 189
$ws.Cells.Item(7, 5).Value = 657

# Row 8
$ws.Cells.Item(8, 2).Value = 44018
$ws.Cells.Item(8, 3).Value = 17152
$ws.Cells.Item(8, 4).Value = 593
$ws.Cells.Item(8, 5).Value = 1650
$ws.Cells.Item(8, 6).Value = 22
$ws.Cells.Item(8, 7).Value = 14.14
$ws.Cells.Item(8, 8).Value = 4.14
$ws.Cells.Item(8, 11).Value = 11675
$ws.Cells.Item(8, 12).Value = 554

# Row 12
$ws.Cells.Item(12, 2).Value = 44018
$ws.Cells.Item(12, 3).Value = 13507
$ws.Cells.Item(12, 4).Value = 515
$ws.Cells.Item(12, 5).Value = 259
$ws.Cells.Item(12, 7).Value = 1.92

# Row 13
$ws.Cells.Item(13, 2).Value = 44018
$ws.Cells.Item(13, 3).Value = 48992
$ws.Cells.Item(13, 4).Value = 1051
$ws.Cells.Item(13, 5).Value = 5897
$ws.Cells.Item(13, 7).Value = 21.03
$ws.Cells.Item(13, 11).Value = 28046

# Row 14
$ws.Cells.Item(14, 2).Value = 44018
$ws.Cells.Item(14, 3).Value = 14407
$ws.Cells.Item(14, 4).Value = 61
$ws.Cells.Item(14, 5).Value = 1842
$ws.Cells.Item(14, 7).Value = 24.03
$ws.Cells.Item(14, 11).Value = 7664

# Row 15
$ws.Cells.Item(15, 2).Value = 44018
$ws.Cells.Item(15, 3).Value = 44375
$ws.Cells.Item(15, 5).Value = 15481
$ws.Cells.Item(15, 7).Value = 46.61
$ws.Cells.Item(15, 11).Value = 33213

# Row 16
$ws.Cells.Item(16, 2).Value = 44017
$ws.Cells.Item(16, 3).Value = 116570
$ws.Cells.Item(16, 4).Value = 3534
$ws.Cells.Item(16, 5).Value = 3170
$ws.Cells.Item(16, 6).Value = 368
$ws.Cells.Item(16, 7).Value = 4.7
$ws.Cells.Item(16, 8).Value = 11.21
$ws.Cells.Item(16, 11).Value = 67498
$ws.Cells.Item(16, 12).Value = 3283

# Row 17
$ws.Cells.Item(17, 2).Value = 44018
$ws.Cells.Item(17, 3).Value = 69904
$ws.Cells.Item(17, 4).Value = 3121
$ws.Cells.Item(17, 5).Value = 20043
$ws.Cells.Item(17, 6).Value = 1263
$ws.Cells.Item(17, 7).Value = 28.67
$ws.Cells.Item(17, 8).Value = 40.47
$ws.Cells.Item(17, 11).Value = 57246
$ws.Cells.Item(17, 12).Value = 69882

# Row 18
$ws.Cells.Item(18, 2).Value = 44017
$ws.Cells.Item(18, 3).Value = 31257
$ws.Cells.Item(18, 4).Value = 1114
$ws.Cells.Item(18, 5).Value = 15110
$ws.Cells.Item(18, 6).Value = 558
$ws.Cells.Item(18, 7).Value = 48.34
$ws.Cells.Item(18, 8).Value = 50.09

# Row 19
$ws.Cells.Item(19, 2).Value = 44018
$ws.Cells.Item(19, 3).Value = 87705
$ws.Cells.Item(19, 4).Value = 6754
$ws.Cells.Item(19, 5).Value = 11603
$ws.Cells.Item(19, 11).Value = 38615

# Row 20
$ws.Cells.Item(20, 2).Value = 44018
$ws.Cells.Item(20, 3).Value = 203376
$ws.Cells.Item(20, 4).Value = 3778
$ws.Cells.Item(20, 5).Value = 26511
$ws.Cells.Item(20, 6).Value = 737
$ws.Cells.Item(20, 7).Value = 13.04

# Row 21
$ws.Cells.Item(21, 2).Value = 44018
$ws.Cells.Item(21, 3).Value = 1249
$ws.Cells.Item(21, 7).Value = 0.48

# Row 22
$ws.Cells.Item(22, 2).Value = 44018
$ws.Cells.Item(22, 3).Value = 1251
$ws.Cells.Item(22, 5).Value = 127
$ws.Cells.Item(22, 7).Value = 10.44
$ws.Cells.Item(22, 11).Value = 1217

# Row 23
$ws.Cells.Item(23, 2).Value = 44018
$ws.Cells.Item(23, 3).Value = 34257
$ws.Cells.Item(23, 4).Value = 1691
$ws.Cells.Item(23, 6).Value = 110
$ws.Cells.Item(23, 7).Value = 6.42
$ws.Cells.Item(23, 8).Value = 6.76
$ws.Cells.Item(23, 11).Value = 28159
$ws.Cells.Item(23, 12).Value = 1628

# Row 24
$ws.Cells.Item(24, 2).Value = 44018
$ws.Cells.Item(24, 3).Value = 20046
$ws.Cells.Item(24, 4).Value = 283
$ws.Cells.Item(24, 5).Value = 1184
$ws.Cells.Item(24, 7).Value = 7.65
$ws.Cells.Item(24, 11).Value = 15470

# Row 25
$ws.Cells.Item(25, 2).Value = 44018
$ws.Cells.Item(25, 3).Value = 66089
$ws.Cells.Item(25, 4).Value = 5897
$ws.Cells.Item(25, 5).Value = 19987
$ws.Cells.Item(25, 7).Value = 30.24
$ws.Cells.Item(25, 8).Value = 39.99

# Row 26
$ws.Cells.Item(26, 2).Value = 44017
$ws.Cells.Item(26, 3).Value = 271684
$ws.Cells.Item(26, 4).Value = 6300
$ws.Cells.Item(26, 5).Value = 7693
$ws.Cells.Item(26, 6).Value = 575
$ws.Cells.Item(26, 8).Value = 9.199999999999999
$ws.Cells.Item(26, 11).Value = 177012
$ws.Cells.Item(26, 12).Value = 6227

# Row 27
$ws.Cells.Item(27, 2).Value = 44018
$ws.Cells.Item(27, 3).Value = 48331
$ws.Cells.Item(27, 4).Value = 2505
$ws.Cells.Item(27, 5).Value = 5798
$ws.Cells.Item(27, 6).Value = 363
$ws.Cells.Item(27, 7).Value = 12
$ws.Cells.Item(27, 8).Value = 14.49

# Row 28
$ws.Cells.Item(28, 2).Value = 44018
$ws.Cells.Item(28, 3).Value = 1166
$ws.Cells.Item(28, 7).Value = 1.23
$ws.Cells.Item(28, 11).Value = 2026

# Row 29
$ws.Cells.Item(29, 2).Value = 44018
$ws.Cells.Item(29, 3).Value = 32061
$ws.Cells.Item(29, 5).Value = 5487
$ws.Cells.Item(29, 7).Value = 19.04
$ws.Cells.Item(29, 11).Value = 28819

# Row 30
$ws.Cells.Item(30, 2).Value = 44018
$ws.Cells.Item(30, 3).Value = 97064
$ws.Cells.Item(30, 4).Value = 2878
$ws.Cells.Item(30, 5).Value = 26887
$ws.Cells.Item(30, 6).Value = 1357
$ws.Cells.Item(30, 7).Value = 27.7
$ws.Cells.Item(30, 8).Value = 47.15

# Row 33
$ws.Cells.Item(33, 2).Value = 44018
$ws.Cells.Item(33, 3).Value = 12293
$ws.Cells.Item(33, 5).Value = 3148
$ws.Cells.Item(33, 7).Value = 25.61

# Row 34
$ws.Cells.Item(34, 2).Value = 44018
$ws.Cells.Item(34, 3).Value = 3423
$ws.Cells.Item(34, 4).Value = 109
$ws.Cells.Item(34, 7).Value = 26.52
$ws.Cells.Item(34, 11).Value = 3028

# Row 35
$ws.Cells.Item(35, 2).Value = 44018
$ws.Cells.Item(35, 3).Value = 74529
$ws.Cells.Item(35, 4).Value = 1398
$ws.Cells.Item(35, 5).Value = 11900
$ws.Cells.Item(35, 6).Value = 446
$ws.Cells.Item(35, 7).Value = 23.78
$ws.Cells.Item(35, 8).Value = 33.06
$ws.Cells.Item(35, 11).Value = 50048
$ws.Cells.Item(35, 12).Value = 1349

# Row 36
$ws.Cells.Item(36, 2).Value = 44018
$ws.Cells.Item(36, 3).Value = 147865
$ws.Cells.Item(36, 4).Value = 7026
$ws.Cells.Item(36, 5).Value = 24783
$ws.Cells.Item(36, 7).Value = 16.76
$ws.Cells.Item(36, 8).Value = 27.87

# Row 37
$ws.Cells.Item(37, 2).Value = 44018
$ws.Cells.Item(37, 3).Value = 8052
$ws.Cells.Item(37, 4).Value = 94
$ws.Cells.Item(37, 5).Value = 120
$ws.Cells.Item(37, 7).Value = 1.49
$ws.Cells.Item(37, 8).Value = 1.06

# Row 38
$ws.Cells.Item(38, 2).Value = 44018
$ws.Cells.Item(38, 3).Value = 38569
$ws.Cells.Item(38, 4).Value = 1474
$ws.Cells.Item(38, 5).Value = 7928
$ws.Cells.Item(38, 7).Value = 20.56
$ws.Cells.Item(38, 8).Value = 8.48

# Row 39
$ws.Cells.Item(39, 2).Value = 44018
$ws.Cells.Item(39, 3).Value = 110137
$ws.Cells.Item(39, 4).Value = 8198
$ws.Cells.Item(39, 5).Value = 10370
$ws.Cells.Item(39, 6).Value = 671

# Row 40
$ws.Cells.Item(40, 2).Value = 44018
$ws.Cells.Item(40, 3).Value = 12436
$ws.Cells.Item(40, 4).Value = 356
$ws.Cells.Item(40, 5).Value = 3592
$ws.Cells.Item(40, 7).Value = 31.12
$ws.Cells.Item(40, 8).Value = 40.17
$ws.Cells.Item(40, 11).Value = 11541
$ws.Cells.Item(40, 12).Value = 356
